$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.938.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.31%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.543.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.78%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.17%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'305.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.11%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'102.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +8.08%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.82%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.09%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.549"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.20%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'37.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.81%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0822"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.05%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'7.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.34%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.19%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.933.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.68%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.587.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.45%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'15.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +6.71%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.876"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.55%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.972.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.19%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.61%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.71%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D22").Value = "'71.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.10%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'253.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.24%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.41%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.55%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'27.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -5.81%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.07%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +9.82%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'10.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.58%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'39.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.96%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.90%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'158.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'2.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.14%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0800"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.30%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'WEMIXToken"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'2.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.00%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'LidoDAOToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'3.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.60%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'18.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.70%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.116"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.75%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'EnergySwap"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'24.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +5.25%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Stellar"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.120"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.31%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.23%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.15%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -8.25%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.04%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.066.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.71%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.02%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'86.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.98%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'9.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.95%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.790.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.61%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.30%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'ordi"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'73.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.68%  "
$ws.Range("E51").Style = "Normal"
